# Weekly refresh of fruit/vegetable prices: reassign Fecha / Volumen /
# Precio minimo / Precio maximo / Precio promedio ponderado / Precio $/Kg
# across the existing data rows (2-16, sheet "Sheet1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target row -> new values for columns D, M, N, O, P, S
$updates = @(
    @{ Row = 2;  D = 44875; M = 400; N = 7000;  O = 7500;  P = 7250;  S = 3625 },
    @{ Row = 3;  D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 },
    @{ Row = 4;  D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 },
    @{ Row = 5;  D = 44882; M = 440; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 6;  D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 },
    @{ Row = 7;  D = 44818; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 8;  D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 },
    @{ Row = 9;  D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 10; D = 44819; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 11; D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 },
    @{ Row = 13; D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 14; D = 44874; M = 300; N = 7500;  O = 8000;  P = 7750;  S = 3875 },
    @{ Row = 15; D = 44881; M = 440; N = 6000;  O = 7000;  P = 6500;  S = 3250 },
    @{ Row = 16; D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $u.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $u.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $u.S   # S: Precio $/Kg
}
